$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(42613.761874999997, 6, 53, 42, 53, 40, 22576, 11398, 1672, 296, 238, 6, 4, "Noun"),
    @(42613.891886574071, 10, 54, 41, 54, 25, 11882, 10246, 1525, 271, 208, 6, 2, "Noun"),
    @(42614.889236111114, 12, 59, 38, 59, 45, 26039, 17393, 2272, 396, 258, 6, 5, "Noun"),
    @(42615.887800925928, 20, 55, 39, 55, 25, 20562, 12723, 2003, 328, 234, 12, 4, "Noun")
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

Write-Output "done"
